$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.516.37"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "1.870.54"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  -2.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.50"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("E6").Value = "  -1.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5082"
$ws.Range("E7").Value = "  -1.84%  "
$ws.Range("E8").Value = "  -2.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08354"
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.13"
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.106"
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.201"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").Value = "1.865.69"
$ws.Range("E13").Value = "  +2.46%  "
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.247"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.010"
$ws.Range("E16").Value = "  -2.28%  "
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.29"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06732"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.68"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.008"
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.905"
$ws.Range("E22").Value = "  -1.92%  "
$ws.Range("D23").Value = "28.555.94"
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("E24").Value = "  -1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.201"
$ws.Range("E25").Value = "  -4.65%  "
$ws.Range("D26").Value = "2.083.03"
$ws.Range("E26").Value = "  +2.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.86"
$ws.Range("E27").Value = "  -3.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.56"
$ws.Range("E28").Value = "  -1.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.414"
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.02"
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("E31").Value = "  -1.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.039"
$ws.Range("E32").Value = "  -1.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.750"
$ws.Range("E33").Value = "  -2.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.610"
$ws.Range("E34").Value = "  -1.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02450"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06595"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.950"
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.039"
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.182"
$ws.Range("E40").Value = "  -1.70%  "
$ws.Range("E41").Value = "  -4.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6355"
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("E44").Value = "  -1.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5997"
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.682"
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.000"
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.212"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.18"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("E51").Value = "  -9.31%  "
